$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 8, leaving only the header row and row 2
$ws.Range("A3:E8").EntireRow.Delete()

# Update row 2 with the new values
$ws.Range("A2").Value = "5Q1"
$ws.Range("B2").Value = "LEKHA"
$ws.Range("C2").Value = "22:50"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-03-18"
$ws.Range("E2").Value = "OOPS"
